$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.427.03'
Set-TextCell $ws 'D3' '1.851.27'
Set-TextCell $ws 'E3' '  +0.20%  '
Set-TextCell $ws 'D4' '1.002'
Set-TextCell $ws 'E4' '  +0.21%  '
Set-TextCell $ws 'D5' '240.72'
Set-TextCell $ws 'E5' '  +0.10%  '
Set-TextCell $ws 'D6' '0.6296'
Set-TextCell $ws 'E6' '  +0.17%  '
Set-TextCell $ws 'D7' '1.002'
Set-TextCell $ws 'E7' '  +0.12%  '
Set-TextCell $ws 'D8' '0.07651'
Set-TextCell $ws 'E8' '  +0.58%  '
Set-TextCell $ws 'D9' '0.2912'
Set-TextCell $ws 'E9' '  -0.59%  '
Set-TextCell $ws 'D10' '24.82'
Set-TextCell $ws 'E10' '  +1.36%  '
Set-TextCell $ws 'D11' '2.166.30'
Set-TextCell $ws 'E11' '  +17.16%  '
Set-TextCell $ws 'D12' '0.07752'
Set-TextCell $ws 'E12' '  +0.12%  '
Set-TextCell $ws 'D13' '5.044'
Set-TextCell $ws 'E13' '  +0.79%  '
Set-TextCell $ws 'D14' '0.6819'
Set-TextCell $ws 'E14' '  +0.33%  '
Set-TextCell $ws 'D15' '0.00001074'
Set-TextCell $ws 'E15' '  -0.85%  '
Set-TextCell $ws 'D16' '83.46'
Set-TextCell $ws 'E16' '  -0.39%  '
Set-TextCell $ws 'D17' '6.181'
Set-TextCell $ws 'E17' '  +0.10%  '
Set-TextCell $ws 'D18' '29.521.68'
Set-TextCell $ws 'E18' '  +0.33%  '
Set-TextCell $ws 'D19' '228.70'
Set-TextCell $ws 'E19' '  -0.01%  '
Set-TextCell $ws 'D20' '12.34'
Set-TextCell $ws 'E20' '  -0.85%  '
Set-TextCell $ws 'D21' '1.002'
Set-TextCell $ws 'E21' '  +0.13%  '
Set-TextCell $ws 'D22' '7.469'
Set-TextCell $ws 'E22' '  +0.03%  '
Set-TextCell $ws 'D23' '1.002'
Set-TextCell $ws 'E23' '  +0.11%  '
Set-TextCell $ws 'D24' '158.09'
Set-TextCell $ws 'E24' '  +0.43%  '
Set-TextCell $ws 'D25' '0.1384'
Set-TextCell $ws 'D26' '8.434'
Set-TextCell $ws 'E26' '  +0.95%  '
Set-TextCell $ws 'D27' '17.71'
Set-TextCell $ws 'E27' '  +0.40%  '
Set-TextCell $ws 'D28' '1.387'
Set-TextCell $ws 'E28' '  +6.54%  '
Set-TextCell $ws 'D29' '1.467'
Set-TextCell $ws 'E29' '  +0.24%  '
Set-TextCell $ws 'D30' '0.05607'
Set-TextCell $ws 'E30' '  +0.48%  '
Set-TextCell $ws 'D31' '4.135'
Set-TextCell $ws 'E31' '  +0.80%  '
Set-TextCell $ws 'D32' '4.061'
Set-TextCell $ws 'E32' '  +0.73%  '
Set-TextCell $ws 'D33' '1.843'
Set-TextCell $ws 'E33' '  +0.06%  '
Set-TextCell $ws 'D35' '0.6966'
Set-TextCell $ws 'E35' '  -1.86%  '
Set-TextCell $ws 'D36' '2.590'
Set-TextCell $ws 'E36' '  +0.17%  '
Set-TextCell $ws 'D37' '0.01806'
Set-TextCell $ws 'E37' '  +0.31%  '
Set-TextCell $ws 'D38' '1.228.15'
Set-TextCell $ws 'E38' '  -0.39%  '
Set-TextCell $ws 'D39' '2.728'
Set-TextCell $ws 'E39' '  -1.57%  '
Set-TextCell $ws 'D40' '6.440'
Set-TextCell $ws 'E40' '  +0.21%  '
Set-TextCell $ws 'D41' '0.9067'
Set-TextCell $ws 'E41' '  +0.09%  '
Set-TextCell $ws 'E42' '  +0.10%  '
Set-TextCell $ws 'D43' '102.46'
Set-TextCell $ws 'E43' '  +0.63%  '
Set-TextCell $ws 'D44' '66.11'
Set-TextCell $ws 'E44' '  +0.04%  '
Set-TextCell $ws 'B45' 'Aptos'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D45' '7.203'
Set-TextCell $ws 'E45' '  +0.04%  '
Set-TextCell $ws 'B46' 'BabyDogeCoin'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 'D46' '0.00000000119'
Set-TextCell $ws 'E46' '  -2.06%  '
Set-TextCell $ws 'D47' '0.4025'
Set-TextCell $ws 'E47' '  +0.10%  '
Set-TextCell $ws 'B48' 'EnergySwap'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws 'D48' '9.039'
Set-TextCell $ws 'E48' '  +0.92%  '
Set-TextCell $ws 'B49' 'Algorand'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell $ws 'D49' '0.1153'
Set-TextCell $ws 'E49' '  +2.87%  '
Set-TextCell $ws 'D50' '1.683'
Set-TextCell $ws 'E50' '  +0.43%  '
Set-TextCell $ws 'D51' '0.05705'
Set-TextCell $ws 'E51' '  -0.07%  '
